$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Name" column header to "Tag" (A1), to support the new
# plasmid-feature columns (resistance genes / origins of replication)
# that get added alongside the existing tag/reagent columns.
$ws.Range("A1").Value = "Tag"

# Re-apply the (unchanged) "Normal" cell style to the header + the data
# row so both carry an explicit style record, matching the refreshed
# config/reagent-naming formatting pass.
$ws.Range("A1:A2").Style = "Normal"

# Leave the cursor on the first data row, under the frozen header.
$ws.Range("A2").Select()
